# Update "想去人数" (F column) counts that increased between two scrapes.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F13").Value = 2375
$ws1.Range("F15").Value = 39
$ws1.Range("F17").Value = 545
$ws1.Range("F18").Value = 162
$ws1.Range("F20").Value = 47
$ws1.Range("F22").Value = 1871
$ws1.Range("F23").Value = 4002
$ws1.Range("F26").Value = 1184
$ws1.Range("F27").Value = 228
$ws1.Range("F28").Value = 2090
$ws1.Range("F32").Value = 112
$ws1.Range("F36").Value = 692
$ws1.Range("F38").Value = 414

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 33

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 2375
$ws4.Range("F15").Value = 33
$ws4.Range("F16").Value = 39
$ws4.Range("F18").Value = 545
$ws4.Range("F19").Value = 162
$ws4.Range("F21").Value = 47
$ws4.Range("F23").Value = 1871
$ws4.Range("F24").Value = 4002
$ws4.Range("F27").Value = 1184
$ws4.Range("F28").Value = 228
$ws4.Range("F29").Value = 2090
$ws4.Range("F33").Value = 112
$ws4.Range("F37").Value = 692
$ws4.Range("F39").Value = 414

$wb.Save()
